# Add files via upload
#
# The uploaded sample workbook now has a second column: column A keeps the
# repo clone URLs, column B carries the branch ref used for that repo
# ("*/master" / "*/main"). The old, now-unused third row (which only held a
# leftover style) is cleared away so the sheet's used range shrinks back
# down to a tidy 2x2 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A already has the two repo URLs - A1 is untouched, but make sure A2
# still reads correctly after the new column/row are introduced.
$ws.Range("A1").Value = "az0508/openrewriteutility.git"
$ws.Range("A2").Value = "az0508/openrewritepreprocessutlity.git"

# New column B: the branch ref that goes with each repo on the same row.
$ws.Range("B1").Value = "*/master"
$ws.Range("B2").Value = "*/main"

# The old A3 (a bare, styled-only placeholder cell) is no longer part of the
# sheet - clear it so the used range becomes A1:B2 again.
$ws.Range("A3").Clear()

# Match the saved selection/active cell from the workbook.
$ws.Range("B2").Select() | Out-Null
